$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D2:E51 to Text format so numeric-looking strings
# (e.g. "611.43") are not auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.378.47"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.670.52"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "611.43"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "143.48"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "2.668.22"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "27.31"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "3.151.12"
$ws.Range("E15").Value = "  +3.68%  "
$ws.Range("D16").Value = "63.197.44"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "2.678.52"
$ws.Range("E18").Value = "  +3.79%  "
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").Value = "341.63"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "6.89"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "67.11"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").Value = "1.65"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "8.65"
$ws.Range("E27").Value = "  +4.82%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "544.19"
$ws.Range("E29").Value = "  +16.17%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("E33").Value = "  +7.09%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "172.90"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  +13.16%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "19.17"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +9.81%  "
$ws.Range("D41").Value = "177.46"
$ws.Range("E41").Value = "  +12.32%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "22.24"
$ws.Range("E44").Value = "  +4.08%  "
$ws.Range("E45").Value = "  +6.68%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0240"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.0962"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "18.71"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D50").Value = "1.75"
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("E51").Value = "  -0.72%  "

# Restore default (General) formatting now that the text values are set,
# so the cell style matches the original workbook (no explicit style index).
$dataRange.ClearFormats()

